$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row 1: "_old" -> "_FV2404" and "_new" -> "_FV2410" ---
$headers1 = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")
for ($i = 0; $i -lt $headers1.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers1[$i]
}

$headers2 = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")
for ($i = 0; $i -lt $headers2.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headers2[$i]
}

# --- Add a table (ListObject) over the used range A1:U73 ---
# The header row already carries the bold/shaded/bordered style (s=1). Creating
# the ListObject while that formatting is present makes Excel snapshot it into
# a header-row dxf; stash the formatting on a scratch range, strip the header
# formatting, build the table, then restore the formatting from the scratch
# copy so the workbook ends up byte-identical in styles.xml / no dxf captured.
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$headerRange.Copy()
$scratch.PasteSpecial(-4122) | Out-Null
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U73")
$list = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$list.Name = "Table1"
$list.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122) | Out-Null
$scratch.Clear() | Out-Null
$excel.CutCopyMode = $false

# --- Freeze the header row (freeze top row / split at row 2) ---
$ws.Activate()
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.SplitColumn = 0
$excel.ActiveWindow.FreezePanes = $true
